$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "filter_or"
$ws.Range("B6").Value = "OwnReality.  Jedem seine Wirklichkeit"
$ws.Range("C6").Value = "OwnReality. À chacun son réel"
$ws.Range("D6").Value = "OwnReality. To Each His Own Reality"

$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "filter_pb"
$ws.Range("B8").Value = "Palais Beauharnais – Vollständiges Inventar der Möbel, Bronzen, Gemälde und anderer Gegenstände"
$ws.Range("C8").Value = "Palais Beauharnais – Inventaire complet des meubles, bronzes, tableaux et autres objets"
$ws.Range("D8").Value = "Palais Beauharnais – full inventory of the furniture, bronzes, paintings and other objects"

$ws.Rows.Item(10).Insert()
$ws.Range("A10").Value = "filter_dfkv"
$ws.Range("B10").Value = "Deutsch-französische Kunstvermittlung von 1870-1961"
$ws.Range("C10").Value = "La réception artistique franco-allemande de 1870-1960 "
$ws.Range("D10").Value = "german and french reception of art between 1870 and 1961"

$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "filter_babue"
$ws.Range("B12").Value = "Bildarchiv Bühler"
$ws.Range("C12").Value = "Archives documentaires et photographiques Hans-Peter Bühler sur la peinture du 19 et 20ème siècles"
$ws.Range("D12").Value = "Picture archive Hans-Peter Bühler on 19th and 20th century painting"

$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = "filter_wikidata"
$ws.Range("B14").Value = "Wikidata"
$ws.Range("C14").Value = "Wikidata"
$ws.Range("D14").Value = "wikidata"

$ws.Rows.Item(21).Insert()
$ws.Range("A21").Value = "filter_av"
$ws.Range("B21").Value = "Architrave – Kunst und Architektur in Paris und Versailles im Spiegel deutscher Reiseberichte des Barock"
$ws.Range("C21").Value = "Architrave – Art et architecture à Paris et Versailles dans les récits de voyageurs allemands à l’époque baroque"
$ws.Range("D21").Value = "Architrave – arts and architecture in Paris and Versailles in accounts by Baroque-Era German travellers"

$ws.Rows.Item(23).Insert()
$ws.Range("A23").Value = "filter_fs"
$ws.Range("B23").Value = "Briefwechsel zwischen Henri Fantin-Latour und Otto Scholderer, 1858–1903"
$ws.Range("C23").Value = "Correspondence between Henri Fantin-Latour and Otto Scholderer, 1858-1903"
$ws.Range("D23").Value = "Correspondence between Henri Fantin-Latour and Otto Scholderer, 1858-1903"

$ws.Rows.Item(25).Insert()
$ws.Range("A25").Value = "filter_ar"
$ws.Range("B25").Value = "Kunstsammlung der Académie Royale de Peinture et de Sculpture"
$ws.Range("C25").Value = "La collection d'art de l’Académie royale de peinture et de sculpture"
$ws.Range("D25").Value = "The art collection of the Académie Royale de Peinture et de Sculpture"
